$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2890
$ws.Range("I64").Value = 2835
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 2835
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2587

$ws.Range("H67").Value = 2890
$ws.Range("I67").Value = 2835
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 2835
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -1977

$ws.Range("H88").Value = 2579.8
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2579.8
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2579.8
$ws.Range("N88").Value = -3391.8

$ws.Range("H91").Value = 2579.8
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2579.8
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2579.8
$ws.Range("N91").Value = -5387.8

$ws.Range("H98").Value = 2415.8333
$ws.Range("I98").Value = 2461.3076
$ws.Range("J98").Value = 2297.6
$ws.Range("K98").Value = 2461.3076
$ws.Range("L98").Value = 2297.6
$ws.Range("M98").Value = -963.3076000000001
$ws.Range("N98").Value = -5293.6

$ws.Range("H100").Value = 1640.9166
$ws.Range("I100").Value = 966
$ws.Range("J100").Value = 2315.8333
$ws.Range("K100").Value = 966
$ws.Range("L100").Value = 2315.8333
$ws.Range("M100").Value = -425
$ws.Range("N100").Value = -3397.8333

$ws.Range("H122").Value = 2415.8333
$ws.Range("I122").Value = 2461.3076
$ws.Range("J122").Value = 2297.6
$ws.Range("K122").Value = 7383.9228
$ws.Range("L122").Value = 6892.799999999999
$ws.Range("M122").Value = -4933.9228
$ws.Range("N122").Value = -11792.8

$ws.Range("H129").Value = 1633.0834
$ws.Range("I129").Value = 792
$ws.Range("J129").Value = 2053.625
$ws.Range("K129").Value = 2376
$ws.Range("L129").Value = 6160.875
$ws.Range("M129").Value = 2624
$ws.Range("N129").Value = -16160.875

$ws.Range("H138").Value = 2600.349
$ws.Range("I138").Value = 3712.7
$ws.Range("J138").Value = 2082.9768
$ws.Range("K138").Value = 11138.1
$ws.Range("L138").Value = 6248.930399999999
$ws.Range("M138").Value = -5998.099999999999
$ws.Range("N138").Value = -16528.9304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 8100
$ws.Range("I53").Value = 1200
$ws.Range("J53").Value = 15000
$ws.Range("K53").Value = 1200
$ws.Range("L53").Value = 15000
$ws.Range("M53").Value = -518

$ws.Range("H61").Value = 5802.609
$ws.Range("I61").Value = 6650.8125
$ws.Range("J61").Value = 3863.8572
$ws.Range("K61").Value = 6650.8125
$ws.Range("L61").Value = 3863.8572
$ws.Range("M61").Value = -6438.8125
$ws.Range("N61").Value = -4287.8572

$ws.Range("H74").Value = 1396.0454
$ws.Range("I74").Value = 479.29413
$ws.Range("J74").Value = 4513
$ws.Range("K74").Value = 479.29413
$ws.Range("L74").Value = 4513
$ws.Range("M74").Value = 394.70587

$ws.Range("H77").Value = 1396.0454
$ws.Range("I77").Value = 479.29413
$ws.Range("J77").Value = 4513
$ws.Range("K77").Value = 2396.47065
$ws.Range("L77").Value = 22565
$ws.Range("M77").Value = 1971.52935

$ws.Range("H88").Value = 3940
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 4300
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 4300
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -5112

$ws.Range("H91").Value = 3940
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 4300
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 4300
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -7108

$ws.Range("H132").Value = 2331.25
$ws.Range("I132").Value = 1978.1177
$ws.Range("J132").Value = 4332.3335
$ws.Range("K132").Value = 5934.3531
$ws.Range("L132").Value = 12997.0005
$ws.Range("M132").Value = -3404.3531

$ws.Range("H136").Value = 5802.609
$ws.Range("I136").Value = 6650.8125
$ws.Range("J136").Value = 3863.8572
$ws.Range("K136").Value = 19952.4375
$ws.Range("L136").Value = 11591.5716
$ws.Range("M136").Value = -17402.4375
$ws.Range("N136").Value = -16691.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 880.8333
$ws.Range("I22").Value = 821.75
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 821.75
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -648.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2232.7856
$ws.Range("I31").Value = 1899.2
$ws.Range("J31").Value = 2418.111
$ws.Range("K31").Value = 1899.2
$ws.Range("L31").Value = 2418.111
$ws.Range("M31").Value = -1604.2
$ws.Range("N31").Value = -3008.111

$ws.Range("H34").Value = 2232.7856
$ws.Range("I34").Value = 1899.2
$ws.Range("J34").Value = 2418.111
$ws.Range("K34").Value = 1899.2
$ws.Range("L34").Value = 2418.111
$ws.Range("M34").Value = -1697.2
$ws.Range("N34").Value = -2822.111

$ws.Range("H62").Value = 3333.3333
$ws.Range("I62").Value = 3333.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3333.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2709.3333
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3333.3333
$ws.Range("I65").Value = 3333.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16666.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13546.6665
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 3830.2354
$ws.Range("I122").Value = 2591.9167
$ws.Range("J122").Value = 6802.2
$ws.Range("K122").Value = 7775.750100000001
$ws.Range("L122").Value = 20406.6
$ws.Range("M122").Value = -5325.750100000001

$ws.Range("H134").Value = 2574.2173
$ws.Range("I134").Value = 2252.0952
$ws.Range("J134").Value = 5956.5
$ws.Range("K134").Value = 6756.285600000001
$ws.Range("L134").Value = 17869.5
$ws.Range("M134").Value = -4221.285600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1885
$ws.Range("I64").Value = 1045
$ws.Range("J64").Value = 2725
$ws.Range("K64").Value = 3135
$ws.Range("L64").Value = 8175
$ws.Range("M64").Value = -2865
$ws.Range("N64").Value = -8715

$ws.Range("H67").Value = 1885
$ws.Range("I67").Value = 1045
$ws.Range("J67").Value = 2725
$ws.Range("K67").Value = 3135
$ws.Range("L67").Value = 8175
$ws.Range("M67").Value = -2199
$ws.Range("N67").Value = -10047

$ws.Range("H131").Value = 22630.906
$ws.Range("I131").Value = 707.5
$ws.Range("J131").Value = 25762.822
$ws.Range("K131").Value = 2122.5
$ws.Range("L131").Value = 77288.466
$ws.Range("M131").Value = 2917.5
$ws.Range("N131").Value = -87368.466

$ws.Range("H137").Value = 3414.4
$ws.Range("I137").Value = 1042.4
$ws.Range("J137").Value = 5786.4
$ws.Range("K137").Value = 3127.2
$ws.Range("L137").Value = 17359.2
$ws.Range("M137").Value = 1972.8
$ws.Range("N137").Value = -27559.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 722.8125
$ws.Range("I113").Value = 332.8889
$ws.Range("J113").Value = 1224.1428
$ws.Range("K113").Value = 332.8889
$ws.Range("L113").Value = 1224.1428
$ws.Range("M113").Value = 1837.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2346.647
$ws.Range("I61").Value = 2309.0952
$ws.Range("J61").Value = 2407.3076
$ws.Range("K61").Value = 2309.0952
$ws.Range("L61").Value = 2407.3076
$ws.Range("M61").Value = -2107.0952
$ws.Range("N61").Value = -2811.3076

$ws.Range("H113").Value = 2346.647
$ws.Range("I113").Value = 2309.0952
$ws.Range("J113").Value = 2407.3076
$ws.Range("K113").Value = 2309.0952
$ws.Range("L113").Value = 2407.3076
$ws.Range("M113").Value = -139.0952000000002
$ws.Range("N113").Value = -6747.3076

$ws.Range("H132").Value = 2993.8
$ws.Range("I132").Value = 2009.7
$ws.Range("J132").Value = 3649.8667
$ws.Range("K132").Value = 6029.1
$ws.Range("L132").Value = 10949.6001
$ws.Range("M132").Value = -3499.1
$ws.Range("N132").Value = -16009.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51248

$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156240

$ws.Range("H132").Value = 2348.24
$ws.Range("I132").Value = 1517.579
$ws.Range("J132").Value = 4978.6665
$ws.Range("K132").Value = 4552.737
$ws.Range("L132").Value = 14935.9995
$ws.Range("M132").Value = -2022.737
